$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Delete")

# Update the "responseBodyRequest" (column D) values for the Delete test-case sheet
$ws.Range("D2").Value  = "Phone Number not registered"
$ws.Range("D3").Value  = "Your account has been deleted"
$ws.Range("D4").Value  = "Phone Number not registered"
$ws.Range("D5").Value  = "Phone Number not registered"
$ws.Range("D8").Value  = "Phone Number not registered"
$ws.Range("D9").Value  = "Phone Number not registered"
$ws.Range("D10").Value = "Phone Number not registered"
$ws.Range("D11").Value = "Phone Number not registered"
$ws.Range("D12").Value = "Phone Number not registered"

# Widen column D so the new, longer text fits (matches author's resize)
$ws.Columns.Item(4).ColumnWidth = 25.830729166666668
